$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KnownIssues")

# D4: reduce exercise number from 3 to 2
$ws.Range("D4").Value = 2

# Row 8: clear Type / Description / Severity, set Sub Type + Description(your own) to the new text
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "there is no Known issue"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("K8").Value = "there is no Known issue"

# Move the active selection to K8
$ws.Range("K8").Select()
